$d = $word.ActiveDocument

# Locate the end of the "Current Coursework" bullet (ends with
# "Introduction to Data Modeling II") so we can insert the new
# "Dean's List" bullet immediately after it.
$search = $d.Content
$found = $search.Find.Execute("Introduction to Data Modeling II", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)

# Re-fetch a clean Range object at the end of the match; reusing the
# Find-affected range object directly can leave InsertXML addressing a
# stale position.
$target = $d.Range($search.End, $search.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="36"/>
              </w:numPr>
              <w:tabs>
                <w:tab w:val="right" w:pos="10080"/>
              </w:tabs>
              <w:ind w:left="-180" w:hanging="270"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:i/>
                <w:iCs/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>Dean&#8217;s List</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
                <w:color w:val="000000" w:themeColor="text1"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>: Spring 2020, Spring 2021</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)
